$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88; this shifts the existing rows 88-230 down
# to become rows 89-231 (matching the target diff, which shows every row's
# data shifted down by one and a brand-new record inserted at row 88).
$ws.Rows("88:88").Insert()

# Populate the newly inserted row 88 with the new weekly record.
$ws.Cells.Item(88, 1).Value = 8
$ws.Cells.Item(88, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(88, 3).Value = "Coquimbo"
$ws.Cells.Item(88, 4).Value = 44797
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(88, 6).Value = 100112037
$ws.Cells.Item(88, 7).Value = "Cebollín"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 800
$ws.Cells.Item(88, 11).Value = 1400
$ws.Cells.Item(88, 12).Value = 1600
$ws.Cells.Item(88, 13).Value = 1500
$ws.Cells.Item(88, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(88, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(88, 16).Value = 250
$ws.Cells.Item(88, 17).Value = 6
$ws.Cells.Item(88, 18).Value = "Hortaliza"
